$wb = $excel.ActiveWorkbook

# --- Sheets before this edit: [1]="总计", [2]="2021-Q4" ---

# 1. Duplicate the existing "2021-Q4" sheet, placing the copy right after it.
#    The copy keeps the old 2021-Q4 data (it becomes the new sheetId=3 tab),
#    while the original sheet2 will be repurposed below to hold the new
#    2022-Q4 numbers. Rename the original out of the way first so the two
#    sheets never collide on the "2021-Q4" name.
$oldQ4 = $wb.Worksheets.Item(2)
$oldQ4.Copy($null, $oldQ4)
$oldQ4.Name = "2022-Q4"

# The newly inserted copy lands immediately after $oldQ4, i.e. position 3.
$newQ4Copy = $wb.Worksheets.Item(3)
$newQ4Copy.Name = "2021-Q4"

# 2. Repurpose the original sheet (still position 2, now named "2022-Q4")
#    and replace its fund table with the new quarter's data.
$ws2022 = $oldQ4

$ws2022.Cells.Item(1, 2).Value = "基金代码"
$ws2022.Cells.Item(1, 3).Value = "基金名称"
$ws2022.Cells.Item(1, 4).Value = "基金规模"
$ws2022.Cells.Item(1, 5).Value = "股票总仓位"
$ws2022.Cells.Item(1, 6).Value = "仓位占比"
$ws2022.Cells.Item(1, 7).Value = "持有市值(亿元)"
$ws2022.Cells.Item(1, 8).Value = "仓位排名"

# Columns B, D, E, F, G hold numeric-looking text in the source data (fund
# code / ratios kept as strings, not numbers). Flip to text format before
# writing so the engine doesn't auto-coerce them to numbers, then clear the
# formatting back off afterwards so no stray style index is left behind.
$ws2022.Range("B2:B3").NumberFormat = "@"
$ws2022.Range("D2:G3").NumberFormat = "@"

$ws2022.Cells.Item(2, 1).Value = 0
$ws2022.Cells.Item(2, 2).Value = "167703"
$ws2022.Cells.Item(2, 3).Value = "德邦量化优选股票（LOF）C"
$ws2022.Cells.Item(2, 4).Value = "0.56"
$ws2022.Cells.Item(2, 5).Value = "88.52"
$ws2022.Cells.Item(2, 6).Value = "1.16"
$ws2022.Cells.Item(2, 7).Value = "0.0065"
$ws2022.Cells.Item(2, 8).Value = 3

$ws2022.Cells.Item(3, 1).Value = 1
$ws2022.Cells.Item(3, 2).Value = "167702"
$ws2022.Cells.Item(3, 3).Value = "德邦量化优选股票（LOF）A"
$ws2022.Cells.Item(3, 4).Value = "0.32"
$ws2022.Cells.Item(3, 5).Value = "88.52"
$ws2022.Cells.Item(3, 6).Value = "1.16"
$ws2022.Cells.Item(3, 7).Value = "0.0037"
$ws2022.Cells.Item(3, 8).Value = 3

$ws2022.Range("B2:B3").ClearFormats()
$ws2022.Range("D2:G3").ClearFormats()

# Row 3's A cell (fund rank index) should carry the same style as row 2's,
# matching the rest of the column.
$ws2022.Cells.Item(2, 1).Copy()
$ws2022.Cells.Item(3, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3. Update the "总计" summary sheet: insert a new row 2 for 2022-Q4 so the
#    existing 2021-Q4 row slides down to row 3 (its data stays the same,
#    only its position index A3 changes from 0 to 1).
$summary = $wb.Worksheets.Item(1)
$summary.Rows.Item(2).Insert()

# The inserted blank row inherits stray formatting on B2:D2 (it has none in
# the final sheet) - strip it back to the default, then restore column A's
# rank-style formatting by copying it down from row 3 (which still carries
# the original style).
$summary.Range("B2:D2").ClearFormats()

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 2
$summary.Cells.Item(2, 4).Value = 0.01

$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(3, 2).Value = "2021-Q4"
$summary.Cells.Item(3, 3).Value = 1
$summary.Cells.Item(3, 4).Value = 0.86

$summary.Cells.Item(3, 1).Copy()
$summary.Cells.Item(2, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Leave the original active tab ("总计") selected, same as before the edit.
$summary.Activate()
